$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the threshold values (B2, B3, C5)
$ws.Range("B2").Value = 5.5
$ws.Range("B3").Value = 5.5
$ws.Range("C5").Value = 25

# Update the active selection to match the saved UI state
$ws.Range("B5:C5").Select()
